# Update database and change read_price algorithm
# - Shift yearly columns (D:H) left by one fiscal year (drop 1396/12, add 1401/12)
# - Shift "publish date" row accordingly, adding the new 1402-01-29 publish dates
# - Shift all financial data rows (11-27) left by one column, adding the new
#   rightmost-column (1401/12) values
#
# NOTE on ordering: brand-new strings that did not exist anywhere in the
# workbook before (the new "1401/12" header, and the two new "1402-01-29"
# publish-date strings) are appended to the shared string table in the order
# they are first written. To reproduce the exact shared-string layout of the
# target workbook, those three assignments (H8, G9, H9) are intentionally
# done last, in that order, after every other (already-existing-string /
# numeric) cell has been updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: fiscal-year-end headers (D:G reuse existing strings) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"

# --- Row 9: publish dates (D:F reuse existing strings) ---
$ws.Range("D9").Value = "1399-03-24 (10)"
$ws.Range("E9").Value = "1400-04-02 (9)"
$ws.Range("F9").Value = "1401-04-08 (9)"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 708263
$ws.Range("E11").Value = 763586
$ws.Range("F11").Value = 1220478
$ws.Range("G11").Value = 1728181
$ws.Range("H11").Value = 1163255

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -439333
$ws.Range("E12").Value = -395344
$ws.Range("F12").Value = -410832
$ws.Range("G12").Value = -630616
$ws.Range("H12").Value = -620833

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 268929
$ws.Range("E13").Value = 368242
$ws.Range("F13").Value = 743965
$ws.Range("G13").Value = 1097565
$ws.Range("H13").Value = 542422

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -3855
$ws.Range("E14").Value = -3408
$ws.Range("F14").Value = -14033
$ws.Range("G14").Value = -28704
$ws.Range("H14").Value = -28196

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (stays "-" across all years) ---
# (no change needed; already "-" in every column)

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = "-"
# E16:H16 already "-"

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 265074
$ws.Range("E17").Value = 364834
$ws.Range("F17").Value = 729933
$ws.Range("G17").Value = 1068862
$ws.Range("H17").Value = 514225

# --- Row 18: هزینه های مالی (Financial expenses) ---
# D18 already "-"
$ws.Range("E18").Value = -2070
$ws.Range("F18").Value = -265
$ws.Range("G18").Value = "-"
# H18 already "-"

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = 10417
$ws.Range("E19").Value = 55680
$ws.Range("F19").Value = 24595
$ws.Range("G19").Value = 42376
$ws.Range("H19").Value = 65182

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 275491
$ws.Range("E20").Value = 418444
$ws.Range("F20").Value = 819944
$ws.Range("G20").Value = 1111238
$ws.Range("H20").Value = 579407

# --- Row 21: مالیات (Tax) ---
# D21 already "-"
$ws.Range("E21").Value = -60552
$ws.Range("F21").Value = -63947
$ws.Range("G21").Value = -89924
$ws.Range("H21").Value = "-"

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 275491
$ws.Range("E22").Value = 357892
$ws.Range("F22").Value = 755997
$ws.Range("G22").Value = 1021313
$ws.Range("H22").Value = 579407

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (stays "-") ---
# (no change needed)

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 275491
$ws.Range("E24").Value = 357892
$ws.Range("F24").Value = 755997
$ws.Range("G24").Value = 1021313
$ws.Range("H24").Value = 579407

# --- Row 25: سود هر سهم پس از کسر مالیات (stays 0) ---
# (no change needed)

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 331168
$ws.Range("E26").Value = 342994
$ws.Range("F26").Value = 245464
$ws.Range("G26").Value = 451008
$ws.Range("H26").Value = 545495

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه (stays 0) ---
# (no change needed)

# --- Brand-new shared strings, added last and in this exact order ---
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("G9").Value = "1402-01-29 (8)"

# H9 ("1402-01-29") looks exactly like a plain ISO date, so a normal .Value
# assignment would be auto-converted by Excel into a date serial number.
# Route it through a text formula and "paste values" so it lands as a plain
# shared string (matching the original cell's style) instead of a date.
$ws.Range("H9").Formula = '="1402-01-29"'
$ws.Range("H9").Copy()
$ws.Range("H9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
